$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.935.87"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.272.02"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.85"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.638"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.68"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("E9").Value = "  +6.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +6.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.54"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "27.32"
$ws.Range("E12").Value = "  +14.59%  "
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "2.611.79"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +6.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "2.274.27"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "43.867.85"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  +7.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.80"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.92"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.35"
$ws.Range("E27").Value = "  +26.69%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.86"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("E32").Value = "  -5.83%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0705"
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -3.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +4.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.50"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0259"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000228"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0994"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.57"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("E45").Value = "  -6.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.48"
$ws.Range("E46").Value = "  +8.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.34"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "1.446.57"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("E51").Value = "  +0.95%  "
